# Add the new "ODI Batting Extra" worksheet as the last sheet in the workbook.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# ---- Header row (bold, centered, thin border - matches other sheets' header style) ----
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---- Data rows ----
# Columns: A=MATCH_CODE(text) B=BATTING_POSITION(number) C=NUM_4(text)
#          D=NUM_6(text) E=PERCENT_RUNS_OF_TOTAL(text) F=MAN_OF_MATCH(text)
$rows = @(
    @{ A = "4379"; B = 2;    C = "0";  D = "0"; E = "0.80%";  F = "NO" },
    @{ A = "4537"; B = 2;    C = "1";  D = "1"; E = "8.84%";  F = "NO" },
    @{ A = "4582"; B = 2;    C = "0";  D = "0"; E = "1.81%";  F = "NO" },
    @{ A = "4585"; B = 2;    C = "16"; D = "0"; E = "52.40%"; F = "YES" },
    @{ A = "4588"; B = 2;    C = "0";  D = "0"; E = "5.84%";  F = "NO" },
    @{ A = "4671"; B = $null; C = $null; D = $null; E = $null; F = "NO" },
    @{ A = "4674"; B = 2;    C = "2";  D = "0"; E = "4.39%";  F = "NO" },
    @{ A = "4675"; B = $null; C = $null; D = $null; E = $null; F = "NO" }
)

$rowIndex = 2
foreach ($r in $rows) {
    $ws.Cells.Item($rowIndex, 1).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 1).Value = $r.A

    if ($null -ne $r.B) {
        $ws.Cells.Item($rowIndex, 2).Value = $r.B
    }

    if ($null -ne $r.C) {
        $ws.Cells.Item($rowIndex, 3).NumberFormat = "@"
        $ws.Cells.Item($rowIndex, 3).Value = $r.C
    }

    if ($null -ne $r.D) {
        $ws.Cells.Item($rowIndex, 4).NumberFormat = "@"
        $ws.Cells.Item($rowIndex, 4).Value = $r.D
    }

    if ($null -ne $r.E) {
        $ws.Cells.Item($rowIndex, 5).NumberFormat = "@"
        $ws.Cells.Item($rowIndex, 5).Value = $r.E
    }

    $ws.Cells.Item($rowIndex, 6).Value = $r.F

    $rowIndex = $rowIndex + 1
}

# Keep first sheet active/selected, matching original workbook view state.
$wb.Worksheets.Item(1).Select()
